$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.660.23"
$ws.Range("E2").Value = "  -0.48%  "
$ws.Range("D3").Value = "2.408.43"
$ws.Range("E3").Value = "  -3.94%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "485.15"
$ws.Range("E5").Value = "  -2.36%  "
$ws.Range("D6").Value = "152.52"
$ws.Range("E6").Value = "  -1.05%  "
$ws.Range("D7").Value = "0.996"
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  +16.09%  "
$ws.Range("D9").Value = "2.425.84"
$ws.Range("E9").Value = "  -3.78%  "
$ws.Range("D10").Value = "0.0994"
$ws.Range("E10").Value = "  -0.06%  "
$ws.Range("D11").Value = "5.65"
$ws.Range("E11").Value = "  -2.42%  "
$ws.Range("D12").Value = "0.334"
$ws.Range("E12").Value = "  -0.73%  "
$ws.Range("E13").Value = "  +1.33%  "
$ws.Range("D14").Value = "2.829.45"
$ws.Range("E14").Value = "  -3.49%  "
$ws.Range("D15").Value = "56.904.34"
$ws.Range("E15").Value = "  -0.25%  "
$ws.Range("D16").Value = "20.71"
$ws.Range("E16").Value = "  -4.01%  "
$ws.Range("E17").Value = "  -3.39%  "
$ws.Range("D18").Value = "2.427.15"
$ws.Range("E18").Value = "  -3.73%  "
$ws.Range("D19").Value = "4.73"
$ws.Range("E19").Value = "  +3.53%  "
$ws.Range("D20").Value = "324.15"
$ws.Range("E20").Value = "  -0.57%  "
$ws.Range("D21").Value = "9.95"
$ws.Range("E21").Value = "  -4.30%  "
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").Value = "  -0.06%  "
$ws.Range("D23").Value = "5.92"
$ws.Range("E23").Value = "  -0.42%  "
$ws.Range("D24").Value = "57.90"
$ws.Range("E24").Value = "  -2.09%  "
$ws.Range("D25").Value = "0.407"
$ws.Range("E25").Value = "  -1.34%  "
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("D27").Value = "0.159"
$ws.Range("E27").Value = "  -2.98%  "
$ws.Range("D28").Value = "2.516.50"
$ws.Range("E28").Value = "  -3.68%  "
$ws.Range("D29").Value = "7.27"
$ws.Range("E29").Value = "  -5.67%  "
$ws.Range("D30").Value = "0.0₃0782"
$ws.Range("E30").Value = "  -4.78%  "
$ws.Range("E31").Value = "  +0.10%  "
$ws.Range("D32").Value = "18.58"
$ws.Range("E32").Value = "  +0.63%  "
$ws.Range("D33").Value = "148.56"
$ws.Range("E33").Value = "  -1.95%  "
$ws.Range("E34").Value = "  -1.36%  "
$ws.Range("E35").Value = "  +1.64%  "
$ws.Range("E36").Value = "  -2.59%  "
$ws.Range("D37").Value = "3.70"
$ws.Range("E37").Value = "  -3.05%  "
$ws.Range("D38").Value = "0.844"
$ws.Range("E38").Value = "  -4.15%  "
$ws.Range("D39").Value = "0.102"
$ws.Range("E39").Value = "  +8.90%  "
$ws.Range("D40").Value = "34.08"
$ws.Range("E40").Value = "  -0.96%  "
$ws.Range("D41").Value = "3.52"
$ws.Range("E41").Value = "  -0.29%  "
$ws.Range("D42").Value = "1.37"
$ws.Range("E42").Value = "  -2.68%  "
$ws.Range("E43").Value = "  +0.05%  "
$ws.Range("D44").Value = "0.590"
$ws.Range("E44").Value = "  -4.33%  "
$ws.Range("D45").Value = "266.73"
$ws.Range("E45").Value = "  -1.29%  "
$ws.Range("D46").Value = "0.0531"
$ws.Range("E46").Value = "  -6.65%  "
$ws.Range("D47").Value = "10.22"
$ws.Range("E47").Value = "  -0.04%  "
$ws.Range("D48").Value = "0.0228"
$ws.Range("E48").Value = "  -1.67%  "
$ws.Range("D49").Value = "4.63"
$ws.Range("E49").Value = "  -7.09%  "
$ws.Range("D50").Value = "17.37"
$ws.Range("E50").Value = "  -3.28%  "
$ws.Range("D51").Value = "1.862.75"
$ws.Range("E51").Value = "  -2.47%  "
